$wb = $excel.ActiveWorkbook

# --- Sheet "ond#Bevestiging": insert a new row 2 (new Bevestiging relation),
#     pushing the existing row 2 down to row 3. ---
$wsBevestiging = $wb.Worksheets.Item("ond#Bevestiging")
$wsBevestiging.Rows.Item(2).Insert()

$wsBevestiging.Range("A2").Value = "https://wegenenverkeer.data.vlaanderen.be/ns/onderdeel#Bevestiging"
$wsBevestiging.Range("B2").Value = "dummy_bevestiging_1"
$wsBevestiging.Range("C2").Value = "dummy_zQp"
$wsBevestiging.Range("D2").Value = "https://wegenenverkeer.data.vlaanderen.be/ns/onderdeel#Pictogram"
$wsBevestiging.Range("E2").Value = "dummy_a"
$wsBevestiging.Range("F2").Value = "dummy_okopD"
$wsBevestiging.Range("G2").Value = "https://wegenenverkeer.data.vlaanderen.be/ns/onderdeel#Funderingsmassief"
$wsBevestiging.Range("H2").Value = "dummy_TyBGmXfXC"
$wsBevestiging.Range("I2").Value = "dummy_dY"
$wsBevestiging.Range("J2").Value = "'False"

# --- Sheet "ond#HoortBij": fold row 3's content into row 2 (keeping A2, G2
#     and J2 unchanged), then delete row 3 entirely. ---
$wsHoortBij = $wb.Worksheets.Item("ond#HoortBij")

$wsHoortBij.Range("B2").Value = "dummy_C_-_dummy_hxOTHWe_-_HoortBij"
$wsHoortBij.Range("C2").Value = "OTLMOW"
$wsHoortBij.Range("D2").Value = "https://wegenenverkeer.data.vlaanderen.be/ns/onderdeel#Pictogram"
$wsHoortBij.Range("E2").Value = "dummy_C"
$wsHoortBij.Range("F2").Value = "dummy_Ek"
$wsHoortBij.Range("H2").Value = "dummy_hxOTHWe"
$wsHoortBij.Range("I2").Value = "dummy_GfaE"

$wsHoortBij.Rows.Item(3).Delete()
